$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet "展览" (exhibitions): update F-column "想去人数" totals
# ============================================================
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(5, 6).Value = 187
$ws1.Cells.Item(6, 6).Value = 12
$ws1.Cells.Item(7, 6).Value = 769
$ws1.Cells.Item(8, 6).Value = 69
$ws1.Cells.Item(9, 6).Value = 9909
$ws1.Cells.Item(10, 6).Value = 53
$ws1.Cells.Item(11, 6).Value = 2682
$ws1.Cells.Item(13, 6).Value = 2406
$ws1.Cells.Item(14, 6).Value = 2700
$ws1.Cells.Item(16, 6).Value = 291
$ws1.Cells.Item(17, 6).Value = 2108
$ws1.Cells.Item(19, 6).Value = 86
$ws1.Cells.Item(20, 6).Value = 372
$ws1.Cells.Item(22, 6).Value = 95
$ws1.Cells.Item(25, 6).Value = 175
$ws1.Cells.Item(26, 6).Value = 605
$ws1.Cells.Item(27, 6).Value = 1301
$ws1.Cells.Item(29, 6).Value = 96
$ws1.Cells.Item(30, 6).Value = 126
$ws1.Cells.Item(32, 6).Value = 1761
$ws1.Cells.Item(33, 6).Value = 2869
$ws1.Cells.Item(34, 6).Value = 4
$ws1.Cells.Item(36, 6).Value = 1009
$ws1.Cells.Item(37, 6).Value = 363
$ws1.Cells.Item(38, 6).Value = 3
$ws1.Cells.Item(39, 6).Value = 1280
$ws1.Cells.Item(40, 6).Value = 67
$ws1.Cells.Item(41, 6).Value = 75
$ws1.Cells.Item(42, 6).Value = 58
$ws1.Cells.Item(44, 6).Value = 31

# ============================================================
# Sheet "演出" (performances): a new event was published and
# inserted at row 7 (it is the earliest upcoming one); rows
# 7-14 shift down to 8-15.
# ============================================================
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Rows.Item(7).Insert()

# New row 7: 北京·次元音浪Million Live⏤番音集结
$ws2.Cells.Item(7, 1).Value = 6
$ws2.Cells.Item(7, 2).Value = "2024-09-22"
$ws2.Cells.Item(7, 3).Value = "北京·次元音浪Million Live⏤番音集结"
$ws2.Cells.Item(7, 4).Value = "学清路38号金码大厦B座 北京想象空间"
$ws2.Cells.Item(7, 5).Value = "2024.09.22 13:00-09.22 16:00"
$ws2.Cells.Item(7, 6).Value = 1
$ws2.Cells.Item(7, 7).Value = 88
$ws2.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=90657"
$ws2.Cells.Item(7, 9).Value = "//i1.hdslb.com/bfs/openplatform/202408/Fn9CSOmf1723477511986.jpeg"

# Column A is the plain running index (row number - 2); Insert()
# dragged the old index values down with their rows, so restamp
# rows 8-15 back to the correct sequential numbers.
for ($r = 8; $r -le 15; $r++) {
    $ws2.Cells.Item($r, 1).Value = $r - 1
}

# ============================================================
# Sheet "本地生活" (local life): update F-column "想去人数"
# ============================================================
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 731
$ws3.Cells.Item(3, 6).Value = 963
$ws3.Cells.Item(4, 6).Value = 116
$ws3.Cells.Item(5, 6).Value = 1907

# ============================================================
# Sheet "全部类型" (all types, a static merged snapshot of the
# other three sheets): update the matching F-column values.
# ============================================================
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 731
$ws4.Cells.Item(3, 6).Value = 963
$ws4.Cells.Item(4, 6).Value = 116
$ws4.Cells.Item(9, 6).Value = 187
$ws4.Cells.Item(10, 6).Value = 12
$ws4.Cells.Item(11, 6).Value = 769
$ws4.Cells.Item(12, 6).Value = 69
$ws4.Cells.Item(13, 6).Value = 9909
$ws4.Cells.Item(14, 6).Value = 53
$ws4.Cells.Item(16, 6).Value = 2682
$ws4.Cells.Item(18, 6).Value = 2406
$ws4.Cells.Item(19, 6).Value = 2700
$ws4.Cells.Item(20, 6).Value = 291
$ws4.Cells.Item(21, 6).Value = 2108
$ws4.Cells.Item(23, 6).Value = 86
$ws4.Cells.Item(24, 6).Value = 372
$ws4.Cells.Item(28, 6).Value = 175
$ws4.Cells.Item(29, 6).Value = 605
$ws4.Cells.Item(30, 6).Value = 1301
$ws4.Cells.Item(32, 6).Value = 126
$ws4.Cells.Item(34, 6).Value = 1762
$ws4.Cells.Item(36, 6).Value = 2869
$ws4.Cells.Item(37, 6).Value = 1009
$ws4.Cells.Item(39, 6).Value = 363
$ws4.Cells.Item(41, 6).Value = 3
$ws4.Cells.Item(44, 6).Value = 1280
$ws4.Cells.Item(45, 6).Value = 67
$ws4.Cells.Item(46, 6).Value = 58
$ws4.Cells.Item(48, 6).Value = 31
